# Weekly data refresh: a new week's record is prepended to the data table
# (row 14), pushing all subsequent records down by one row. The last
# existing record (old row 80) ends up at row 81, and the sheet's used
# range grows from A1:R80 to A1:R81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14, shifting rows 14..80 down to 15..81.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = "2023-05-25"
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = 100112001
$ws.Range("G14").Value = "Berenjena"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 70
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 6429
$ws.Range("N14").Value = "$/caja 60 unidades"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 107
$ws.Range("Q14").Value = 60
$ws.Range("R14").Value = "Hortaliza"
